$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 12.92687315769199
$ws.Range("D2").Value = 6.219344393504453
$ws.Range("E2").Value = 13.15998055045497
$ws.Range("F2").Value = 49.3025153564332
$ws.Range("G2").Value = 62.0259216775334
$ws.Range("H2").Value = 22.63425151667966
$ws.Range("J2").Value = 11.20221831876756
$ws.Range("K2").Value = 20.44810689891228
$ws.Range("L2").Value = 9.170871096079406
$ws.Range("M2").Value = 21.5259513024024
$ws.Range("C3").Value = 12.91598200217053
$ws.Range("D3").Value = 6.206201455712169
$ws.Range("E3").Value = 13.18292983078363
$ws.Range("F3").Value = 49.42237644100456
$ws.Range("G3").Value = 62.16784841183699
$ws.Range("H3").Value = 22.70916845062349
$ws.Range("J3").Value = 11.23125199178867
$ws.Range("K3").Value = 20.15897159095709
$ws.Range("L3").Value = 9.182789173972527
$ws.Range("M3").Value = 21.42713859593472
$ws.Range("C4").Value = 12.91179531136909
$ws.Range("D4").Value = 6.198470491780769
$ws.Range("E4").Value = 13.19851721582488
$ws.Range("F4").Value = 49.50796873855926
$ws.Range("G4").Value = 62.27293307335501
$ws.Range("H4").Value = 22.75960943153644
$ws.Range("J4").Value = 11.25016045828546
$ws.Range("K4").Value = 19.9831549704892
$ws.Range("L4").Value = 9.190525006187491
$ws.Range("M4").Value = 21.36974495790199
$ws.Range("C5").Value = 12.91071947203517
$ws.Range("D5").Value = 6.195406709469755
$ws.Range("E5").Value = 13.20524585126817
$ws.Range("F5").Value = 49.54585357889801
$ws.Range("G5").Value = 62.32024432317533
$ws.Range("H5").Value = 22.78127884525524
$ws.Range("J5").Value = 11.25813837983579
$ws.Range("K5").Value = 19.91202266055362
$ws.Range("L5").Value = 9.193782891605318
$ws.Range("M5").Value = 21.34719732465772
$ws.Range("C6").Value = 12.91057893284833
$ws.Range("D6").Value = 6.194903245520501
$ws.Range("E6").Value = 13.20638589505507
$ws.Range("F6").Value = 49.55232549116852
$ws.Range("G6").Value = 62.32837070574784
$ws.Range("H6").Value = 22.78494425967722
$ws.Range("J6").Value = 11.25947958694979
$ws.Range("K6").Value = 19.90024469238386
$ws.Range("L6").Value = 9.194330241195484
$ws.Range("M6").Value = 21.34350454921993
$ws.Range("C7").Value = 12.911778248579
$ws.Range("D7").Value = 6.198428819812851
$ws.Range("E7").Value = 13.19860643504916
$ws.Range("F7").Value = 49.50846751429773
$ws.Range("G7").Value = 62.27355298783291
$ws.Range("H7").Value = 22.75989716458053
$ws.Range("J7").Value = 11.2502669469549
$ws.Range("K7").Value = 19.98219346263943
$ws.Range("L7").Value = 9.190568515670135
$ws.Range("M7").Value = 21.36943744745719
$ws.Range("C8").Value = 12.92259990829525
$ws.Range("D8").Value = 6.214743251972175
$ws.Range("E8").Value = 13.16758309595676
$ws.Range("F8").Value = 49.34134673926478
$ws.Range("G8").Value = 62.07112094095165
$ws.Range("H8").Value = 22.65915945713436
$ws.Range("J8").Value = 11.21200497219155
$ws.Range("K8").Value = 20.34810757045811
$ws.Range("L8").Value = 9.17489388090352
$ws.Range("M8").Value = 21.49121058198522
$ws.Range("C9").Value = 12.96358148153196
$ws.Range("D9").Value = 6.249360808860578
$ws.Range("E9").Value = 13.11860416741369
$ws.Range("F9").Value = 49.1093031329804
$ws.Range("G9").Value = 61.8174704363006
$ws.Range("H9").Value = 22.4969707875404
$ws.Range("J9").Value = 11.1455300086008
$ws.Range("K9").Value = 21.07558873050194
$ws.Range("L9").Value = 9.14745785995593
$ws.Range("M9").Value = 21.7552475381377
$ws.Range("C10").Value = 13.00560831215182
$ws.Range("D10").Value = 6.276308805014999
$ws.Range("E10").Value = 13.08982687152624
$ws.Range("F10").Value = 48.99777841218268
$ws.Range("G10").Value = 61.71970138779156
$ws.Range("H10").Value = 22.3995215568472
$ws.Range("J10").Value = 11.10187115370012
$ws.Range("K10").Value = 21.61124045422961
$ws.Range("L10").Value = 9.129292104379189
$ws.Range("M10").Value = 21.96353848714879
$ws.Range("C11").Value = 13.02727988555498
$ws.Range("D11").Value = 6.288878128062485
$ws.Range("E11").Value = 13.07829565072881
$ws.Range("F11").Value = 48.95996442322909
$ws.Range("G11").Value = 61.69468372492757
$ws.Range("H11").Value = 22.35993639099033
$ws.Range("J11").Value = 11.08312687788747
$ws.Range("K11").Value = 21.85416727019204
$ws.Range("L11").Value = 9.121455911754744
$ws.Range("M11").Value = 22.06115369238966
$ws.Range("C12").Value = 13.03584984841527
$ws.Range("D12").Value = 6.293680779955782
$ws.Range("E12").Value = 13.07415294935575
$ws.Range("F12").Value = 48.94751140805138
$ws.Range("G12").Value = 61.6880233898639
$ws.Range("H12").Value = 22.34563149659371
$ws.Range("J12").Value = 11.07618887612698
$ws.Range("K12").Value = 21.94596877670428
$ws.Range("L12").Value = 9.118549674370144
$ws.Range("M12").Value = 22.09850766707523
$ws.Range("C13").Value = 13.03398805485141
$ws.Range("D13").Value = 6.29264456188938
$ws.Range("E13").Value = 13.07503520120316
$ws.Range("F13").Value = 48.95011025553925
$ws.Range("G13").Value = 61.68933246796662
$ws.Range("H13").Value = 22.34868179490929
$ws.Range("J13").Value = 11.07767598832539
$ws.Range("K13").Value = 21.92620739976239
$ws.Range("L13").Value = 9.119172869539334
$ws.Range("M13").Value = 22.09044587300383
$ws.Range("C14").Value = 13.02797768121439
$ws.Range("D14").Value = 6.289272394935673
$ws.Range("E14").Value = 13.07795034260901
$ws.Range("F14").Value = 48.95890245640565
$ws.Range("G14").Value = 61.69407931200413
$ws.Range("H14").Value = 22.35874577131867
$ws.Range("J14").Value = 11.08255287942681
$ws.Range("K14").Value = 21.86172409921721
$ws.Range("L14").Value = 9.121215589943683
$ws.Range("M14").Value = 22.06421915642069
$ws.Range("C15").Value = 13.02434336215011
$ws.Range("D15").Value = 6.287212383073129
$ws.Range("E15").Value = 13.07976510104267
$ws.Range("F15").Value = 48.96453122255932
$ws.Range("G15").Value = 61.69735368560326
$ws.Range("H15").Value = 22.36499955703678
$ws.Range("J15").Value = 11.08556094493268
$ws.Range("K15").Value = 21.82219911612492
$ws.Range("L15").Value = 9.122474770518249
$ws.Range("M15").Value = 22.04820454704567
$ws.Range("C16").Value = 13.00424309118377
$ws.Range("D16").Value = 6.275493471597443
$ws.Range("E16").Value = 13.09061181868076
$ws.Range("F16").Value = 49.00051036375465
$ws.Range("G16").Value = 61.7217292476107
$ws.Range("H16").Value = 22.40220426444969
$ws.Range("J16").Value = 11.10311855877096
$ws.Range("K16").Value = 21.59534252327812
$ws.Range("L16").Value = 9.129812793392263
$ws.Range("M16").Value = 21.9572148035905
$ws.Range("C17").Value = 12.9925636481341
$ws.Range("D17").Value = 6.268382612641282
$ws.Range("E17").Value = 13.09766515559965
$ws.Range("F17").Value = 49.02589767752328
$ws.Range("G17").Value = 61.74167803347964
$ws.Range("H17").Value = 22.42624547769027
$ws.Range("J17").Value = 11.11417516638029
$ws.Range("K17").Value = 21.45592321977158
$ws.Range("L17").Value = 9.134423700814475
$ws.Range("M17").Value = 21.90211296306819
$ws.Range("C18").Value = 12.98608652526471
$ws.Range("D18").Value = 6.264321970212659
$ws.Range("E18").Value = 13.10186887902184
$ws.Range("F18").Value = 49.04171534476104
$ws.Range("G18").Value = 61.75498273741326
$ws.Range("H18").Value = 22.44051980176273
$ws.Range("J18").Value = 11.1206397397389
$ws.Range("K18").Value = 21.37566770162899
$ws.Range("L18").Value = 9.137116029924515
$ws.Range("M18").Value = 21.87069066668683
$ws.Range("C19").Value = 12.98393491155243
$ws.Range("D19").Value = 6.262952197219635
$ws.Range("E19").Value = 13.10331741663594
$ws.Range("F19").Value = 49.04727943267544
$ws.Range("G19").Value = 61.75980140520403
$ws.Range("H19").Value = 22.44542943449356
$ws.Range("J19").Value = 11.12284659990213
$ws.Range("K19").Value = 21.34848592073057
$ws.Range("L19").Value = 9.138034530140908
$ws.Range("M19").Value = 21.86009880486158
$ws.Range("C20").Value = 12.99378207130922
$ws.Range("D20").Value = 6.269136552043429
$ws.Range("E20").Value = 13.0968991220379
$ws.Range("F20").Value = 49.02306928480083
$ws.Range("G20").Value = 61.73936486507179
$ws.Range("H20").Value = 22.42364002273125
$ws.Range("J20").Value = 11.11298729719241
$ws.Range("K20").Value = 21.47077199996385
$ws.Range("L20").Value = 9.133928697688553
$ws.Range("M20").Value = 21.90795077933846
$ws.Range("C21").Value = 13.02973324271907
$ws.Range("D21").Value = 6.290261731182436
$ws.Range("E21").Value = 13.07708802072125
$ws.Range("F21").Value = 48.95626925857499
$ws.Range("G21").Value = 61.69260858479817
$ws.Range("H21").Value = 22.35577112048312
$ws.Range("J21").Value = 11.08111607842012
$ws.Range("K21").Value = 21.88067019887665
$ws.Range("L21").Value = 9.120613936117753
$ws.Range("M21").Value = 22.0719121972194
$ws.Range("C22").Value = 13.05534562083867
$ws.Range("D22").Value = 6.304317615251914
$ws.Range("E22").Value = 13.06544532754788
$ws.Range("F22").Value = 48.92349327545059
$ws.Range("G22").Value = 61.67845440679087
$ws.Range("H22").Value = 22.3154097142113
$ws.Range("J22").Value = 11.06121906258328
$ws.Range("K22").Value = 22.14742082644009
$ws.Range("L22").Value = 9.112268304617906
$ws.Range("M22").Value = 22.18132831427339
$ws.Range("C23").Value = 13.04148352022613
$ws.Range("D23").Value = 6.296793489261329
$ws.Range("E23").Value = 13.07153996807157
$ws.Range("F23").Value = 48.93998815879749
$ws.Range("G23").Value = 61.6845032988712
$ws.Range("H23").Value = 22.33658490770363
$ws.Range("J23").Value = 11.07175329243706
$ws.Range("K23").Value = 22.00518204077948
$ws.Range("L23").Value = 9.116690024172676
$ws.Range("M23").Value = 22.12273183303574
$ws.Range("C24").Value = 12.9932304820082
$ws.Range("D24").Value = 6.268795610309348
$ws.Range("E24").Value = 13.09724498271351
$ws.Range("F24").Value = 49.02434419485989
$ws.Range("G24").Value = 61.74040493107358
$ws.Range("H24").Value = 22.4248165389525
$ws.Range("J24").Value = 11.11352399646027
$ws.Range("K24").Value = 21.46405917626998
$ws.Range("L24").Value = 9.134152359445267
$ws.Range("M24").Value = 21.90531070197854
$ws.Range("C25").Value = 12.95039114117475
$ws.Range("D25").Value = 6.239723535204218
$ws.Range("E25").Value = 13.13058695523635
$ws.Range("F25").Value = 49.16176601816591
$ws.Range("G25").Value = 61.87061380086951
$ws.Range("H25").Value = 22.53704559334897
$ws.Range("J25").Value = 11.16260095476497
$ws.Range("K25").Value = 20.87823123443275
$ws.Range("L25").Value = 9.154528764194408
$ws.Range("M25").Value = 21.6812249600051
